# issue #5: add legislator_id, name, date into dataframe
#
# Every other sheet's output in this legislator-property-disclosure workbook
# is keyed by (date, legislator_name, legislator_id). The "股票" (stocks)
# sheet was missing those three trailing columns, so add them:
#   H: date              (2011-11-28, the filing date for this disclosure)
#   I: legislator_name   (林德福)
#   J: legislator_id     (908)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$date = "2011-11-28"
$legislatorName = "林德福"
$legislatorId = 908

# last row currently holding stock data (column A carries the row index)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# --- headers (row 1), matching the bold/centered/bordered style already
#     used by the existing B1:G1 header cells -------------------------------
$ws.Cells.Item(1, 8).Value = "date"
$ws.Cells.Item(1, 9).Value = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

$headerRange = $ws.Range($ws.Cells.Item(1, 8), $ws.Cells.Item(1, 10))
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- data rows --------------------------------------------------------------
# Column H holds a literal "yyyy-mm-dd" string (not an Excel date serial),
# so format it as text before writing the value.
$dateDataRange = $ws.Range($ws.Cells.Item(2, 8), $ws.Cells.Item($lastRow, 8))
$dateDataRange.NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = $date
    $ws.Cells.Item($r, 9).Value = $legislatorName
    $ws.Cells.Item($r, 10).Value = $legislatorId
}
